$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: insert a new "Meta description" paragraph right after the
# title (Heading1) paragraph.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaFragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Discover the stunning visuals and innovative Scatter Pays feature in Aloha Fruit Bonanza. With an RTP of 97.01%, play for free and win up to 7,500x your stake.</w:t></w:r></w:p>'
$metaPara.Range.InsertXML($metaFragment)

# ---------------------------------------------------------------------
# Change 2: remove the duplicated bold title paragraph that used to sit
# right before the closing italic meta-description paragraph.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($count - 1)
$dupTitlePara.Range.Delete()

# ---------------------------------------------------------------------
# Change 3: rewrite the text of the final (italic) paragraph with the
# new image-generation prompt, keeping its italic formatting and the
# leading empty run intact.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$promptFragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Prompt: Create a cartoon-style feature image for the game &quot;Aloha Fruit Bonanza&quot;. The image should feature a happy Maya warrior wearing glasses. The design should have a tropical feel, with bright colors and a beach background. It should also prominently feature fruit symbols from the game, such as watermelon, coconut cocktail, and dragon fruit cocktail, as well as the red number seven. The warrior should be holding a slot machine lever, and there should be cascading symbols falling around him. The overall feeling of the image should be fun and exciting, capturing the lightheartedness of the game.</w:t></w:r></w:p>'
$lastPara.Range.InsertXML($promptFragment)

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
